# Generate Report for Handoff
#
# The "b.md" file has moved from "Handed back: in sync with en-US" to
# "Ready for handoff" with a freshly generated handoff package
# (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf). Update the Overview
# sheet and each locale sheet (zh-cn, de-de) to reflect the new status,
# handoff file name and handoff datetime for the b.md row (row 3), and
# update the matching hyperlink display text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" file.
#   B3 = zh-cn status, C3 = de-de status, D3 = Latest Handoff Date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-29-19 14:29:14"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" file.
#   C3 = Status, D3 = Latest Handoff File, E3 = Latest Handoff Datetime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-19 14:29:11"

foreach ($link in $zhcn.Hyperlinks) {
    if ($link.Range.Address() -eq '$D$3') {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" file.
#   C3 = Status, D3 = Latest Handoff File, E3 = Latest Handoff Datetime
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-19 14:29:14"

foreach ($link in $dede.Hyperlinks) {
    if ($link.Range.Address() -eq '$D$3') {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
